# Append-refresh of the scraper output: the source run that produced this
# workbook was re-executed at 2025-09-23 06:27:21 (JST) and returned a
# smaller, re-ranked result set (22 listings -> 7 listings). Rewrite the
# "ランサーズ" sheet in place to match the new scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-09-23 06:27:21"

# --- Row data for the new result set (rows 2-8) ------------------------
# Columns: A=取得日時 B=タイトル C=カテゴリ D=価格 E=締切 F=URL G=優先度スコア H=スキル概要
$rows = @(
    @{ B = "初回 「AIで笑顔を検出し、2秒クリップを無劣化で自動切り出すWindowsツール開発(予算10万円)」"; C = "システム開発"; D = "100,000 円 ~ 200,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5398662"; G = 413; H = "🔥AI,Ai ◆ツール,開発" },
    @{ B = "【急募】出品・在庫管理ツール開発と保守対応者募集"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5398562"; G = 163; H = "◆ツール,開発 ◇管理" },
    @{ B = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5251319"; G = 135; H = "◆ツール,スクレイピング ◇サイト" },
    @{ B = "EC多プラットフォーム展開在庫・価格連携ツールの作成"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5398432"; G = 80; H = "◆ツール" },
    @{ B = "運送会社の作業予定表 WEBシステムの修正開発についての相談"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5398772"; G = 78; H = "◆開発" },
    @{ B = "【急募】災害時に備えた「ピジョンクラウド」でのシステムづくり、運用サポートの依頼"; C = "システム開発"; D = "100,000 円 ~ 200,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5398657"; G = 33; H = $null },
    @{ B = "【急募】オリジナルシャンパンのシミュレーションページ制作"; C = "システム開発"; D = "100,000 円 ~ 200,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5398736"; G = 18; H = $null }
)

# --- Drop every stale listing (old rows 9-23) and any content below the
# new result set in columns A-H. This also shrinks the sheet's used range
# (dimension) down to the surviving rows.
$ws.Range("A9:H23").Clear()

# --- Remove every existing hyperlink object; they will be re-created
# below, scoped exactly to the surviving F2:F8 cells with fresh URLs.
$ws.Range("A1").Hyperlinks.Delete()

# --- Rewrite rows 2-8 ----------------------------------------------------
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E

    $fcell = $ws.Cells.Item($r, 6)
    $fcell.Value = $data.F
    $ws.Hyperlinks.Add($fcell, $data.F)
    $fcell.Style = "Hyperlink"

    $ws.Cells.Item($r, 7).Value = $data.G

    $hcell = $ws.Cells.Item($r, 8)
    if ($data.H -eq $null) {
        $hcell.ClearContents()
    } else {
        $hcell.Value = $data.H
    }
}

# --- Column H widened from 17 to 19 characters --------------------------
# This host's ColumnWidth setter stores (value + 5/6) into the OOXML
# <col width>, matching how the sheet's existing integer widths read back
# 5/6 lower than their stored value (e.g. stored 17 reads as 16.17). Back
# the offset out so the persisted width lands exactly on 19, consistent
# with the plain-integer widths the rest of the sheet already uses.
$ws.Columns.Item(8).ColumnWidth = 19 - 5/6
